$d = $word.ActiveDocument

# --- Helpers -----------------------------------------------------------
# Grab the raw OOXML for just the first <w:p>...</w:p> that the Range.XML
# call returns (Range.XML tends to tack on a following empty paragraph
# marker, so only take the first element).
function Get-ParaFragment($p) {
    $full = $p.Range.get_XML($false)
    if ($full -match "(?s)(<w:p[ >].*?</w:p>)") {
        return $matches[1]
    }
    return $null
}

# Replace the visible text of a paragraph while leaving every other bit of
# markup (leading empty <w:r/>, rPr on the text run, pPr, etc.) untouched,
# by doing a minimal text-only substitution on the paragraph's own XML and
# re-inserting it with InsertXML (Find/Execute ReplaceAll instead tends to
# silently swallow empty sibling runs).
function Replace-ParaText($oldText, $newText) {
    foreach ($p in $d.Paragraphs) {
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $oldText) {
            $frag = Get-ParaFragment $p
            if ($null -eq $frag) { continue }
            # Drop the w14:paraId / rsid* attributes the XML export adds to
            # <w:p> so the saved markup keeps matching the source's bare <w:p>.
            $frag = $frag -replace '<w:p [^>]*>', '<w:p>'
            $escaped = $newText -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
            $newFrag = $frag -replace '(?s)<w:t([^>]*)>.*?</w:t>', ("<w:t`$1>" + $escaped + '</w:t>')
            $p.Range.InsertXML($newFrag)
            return $true
        }
    }
    return $false
}

# --- Edits ---------------------------------------------------------------

# Title / heading text (appears twice: Heading1 and the bold summary line)
Replace-ParaText "Play Diamond Mystic Free - Scratch-Off Style Slot Game" "Play Diamond Mystic Free - Unique Lotto-Style Slot Game"
Replace-ParaText "Play Diamond Mystic Free - Scratch-Off Style Slot Game" "Play Diamond Mystic Free - Unique Lotto-Style Slot Game"

# "What we like" bullet list
Replace-ParaText "Flexible betting range starting from €0.01 to a maximum of €2.00" "Unique gameplay mechanics resembling a scratch-off ticket or lotto game"
Replace-ParaText "Minimal bets with decent payouts" "Decent payouts with minimal bets, starting from €0.01"
Replace-ParaText "Visually appealing graphics with a simple design" "Flexible betting range from €0.01 to €2.00"
Replace-ParaText "Fast-paced gameplay" "Simple yet visually appealing graphics"

# "What we don't like" bullet list
Replace-ParaText "Lacks additional features in the game's internal menu" "Lack of additional features in the internal menu"
Replace-ParaText "No Return to Player percentage information available" "No inclusion of Return to Player percentage"

# Closing italic summary line
Replace-ParaText "Experience the excitement of Diamond Mystic - a scratch-off style slot game with flexible betting options and minimal bets. Play for free today!" "Play Diamond Mystic for free and experience a unique lotto-style slot game with decent payouts."
